# Assignment6.xlsx edit script
# - Measures sheet (Table1): insert a new "Measure Folder" column (position 4),
#   shift "Measure Description" to column E, rewrite the DAX expressions in
#   column B with multi-line formatting, set "No Folder Defined" in the new
#   column D for every data row, and refresh the Measure Description text.
# - Source Information sheet (table "Source"): bump Table No for row 2 to 2
#   and rewrite the Modification Description text; extend table ref to I3.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: Measures
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# Insert a new column before column D ("Measure Description" -> shifts to E).
$ws1.Columns.Item(4).Insert()

# Column widths: new D = 30 (matches columns B/C), E keeps the old 50 width
# that travelled with the shifted data. ColumnWidth is specified in Excel's
# "characters" unit, which renders ~0.8333 wider in the saved <col width>;
# subtract that padding so the serialized width is exactly 30.
$ws1.Columns.Item(4).ColumnWidth = 29.1666666666667

# Resize the table to include the new column before writing header text, so
# the table machinery doesn't cache a stale column name.
$lo1 = $ws1.ListObjects.Item(1)
$lo1.Resize($ws1.Range("A1:E9"))

# Header row
$ws1.Range("D1").Value = "Measure Folder"
$ws1.Range("E1").Value = "Measure Description"

# Row 2 - SalesMTD
$ws1.Range("B2").Value = @'

TOTALMTD(
    SUM (FinancialData[Sales]), 
    DATESMTD('FinancialData'[Date])
)
'@
$ws1.Range("D2").Value = "No Folder Defined"
$ws1.Range("E2").Value = "This calculation is finding the total Sales amount for the month-to-date (MTD). The TOTALMTD function totals the sum of the Sales column in the FinancialData table, while the DATESMTD function looks at the Date column in the FinancialData table to only consider the sales that have occurred"

# Row 3 - SalesQTD
$ws1.Range("B3").Value = @'

TOTALQTD(
    SUM(FinancialData[Sales]), 
    DATESQTD('FinancialData'[Date])
)
'@
$ws1.Range("D3").Value = "No Folder Defined"
$ws1.Range("E3").Value = "This calculation will return the total amount of sales up to and including the current quarter of the given year in the 'FinancialData' table. The sum of the sales from the FinancialData table will be taken and then filtered by the specific dates from the current quarter. The result of this calculation is the total"

# Row 4 - SalesYTD
$ws1.Range("B4").Value = @'

TOTALYTD ( 
    SUM ( FinancialData[Sales] ), 
    DATESYTD ( 'FinancialData'[Date] ) 
)
'@
$ws1.Range("D4").Value = "No Folder Defined"
$ws1.Range("E4").Value = "This calculation sums the total sales from the current year up to date from the FinancialData table. It takes into account the date and only counts sales from the present year."

# Row 5 - Previous Month Sales MTD
$ws1.Range("B5").Value = @'

CALCULATE(
    SUM(FinancialData[Sales]),
    PARALLELPERIOD(
        DATESMTD(FinancialData[Date].[Date]), 
        -1,
        MONTH
    )
)
'@
$ws1.Range("D5").Value = "No Folder Defined"
$ws1.Range("E5").Value = "This calculation finds the total sales for the current month-to-date (MTD), as well as the sales for the same time period in the previous month. It does this by using the Parallel Period function to change the date context of the total sales being calculated."

# Row 6 - Previous Quarter Sales QTD
$ws1.Range("B6").Value = @'

CALCULATE(
    SUM(FinancialData[Sales]),
    DATESYTD(
        PARALLELPERIOD(
            FinancialData[Date].[Date],
            -1,
            QUARTER
        )
    )
)
'@
$ws1.Range("D6").Value = "No Folder Defined"
$ws1.Range("E6").Value = "This calculation is getting the total sum of sales for the current year, but only for the same quarter as the previous year. For example, if the current day was in the 3rd Quarter, this calculation would get the total sales for the 3rd Quarter of the current year and the 3rd Quarter of the"

# Row 7 - Month(%)
$ws1.Range("B7").Value = @'

(
    [Previous Month Sales MTD] - [SalesMTD]
) / [Previous Month Sales MTD]
'@
$ws1.Range("D7").Value = "No Folder Defined"
$ws1.Range("E7").Value = "This calculation calculates the percentage change in monthly sales from one month to the previous month. It is found by subtracting the current month's total sales from the previous month's total sales and then dividing by the previous month's total sales."

# Row 8 - Previous Year Current Month Sales
$ws1.Range("B8").Value = @'

CALCULATE (
    SUM ( FinancialData[Sales] ),
    DATESYTD ( SAMEPERIODLASTYEAR ( FinancialData[Date].[Date] ) )
)
'@
$ws1.Range("D8").Value = "No Folder Defined"
$ws1.Range("E8").Value = "This calculation is summing the total sales for the current year to date, compared to the same time period (year to date) in the previous year. The calculation is useful for observing year-over-year growth in sales."

# Row 9 - First Half Year Sales
$ws1.Range("B9").Value = @'

CALCULATE (
   SUM (FinancialData[Sales]),
   DATESBETWEEN (FinancialData[Date],
                  DATE (2014, 1, 1),
                  DATE (2014, 6, 31))
)
'@
$ws1.Range("D9").Value = "No Folder Defined"
$ws1.Range("E9").Value = "This calculation is used to find the sum of sales from a financial data table from January 1st, 2014 to June 30th, 2014. It uses the calculate and sum functions to calculate the total sales value from those dates."

# ---------------------------------------------------------------------------
# Sheet 2: Source Information
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

$ws2.Range("A2").Value = 2

$ws2.Range("I2").Value = @'
1. "Changed Type" is a command to change the data type of each column in the FinancialData_Table to the specified type.
2. The table's "COGS" column has been renamed to "Cost".
3. This means the types of columns Profit, Cost, Sales, Gross Sales, Discounts, Date and Units Sold are changed to Currency, Date and Integer types respectively.
4. "Added Custom" adds a new column to the table called "Total Cost" which is the sum of the two existing columns "Cost" and "Discounts".
5. This sentence changes the data type of the column "Total Cost" to the currency type in the table "Added Custom."

'@

# Extend the "Source" table by one row (matches the committed table ref
# A1:I3) even though no extra row of data was populated in the source edit.
$lo2 = $ws2.ListObjects.Item(1)
$lo2.Resize($ws2.Range("A1:I3"))
